$d = $word.ActiveDocument
$r = $d.Range(1, 7)
$f = $r.Find
$f.ClearFormatting()
$f.Replacement.ClearFormatting()
$res = $f.Execute("ANNEX ", $false, $false, $false, $false, $false, $true, 1, $false, "ANNEX ", 2)
Write-Host "result: $res"
